$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply edits in the same order the strings appear in the authored
# sharedStrings table so newly-introduced shared strings land in the
# expected order: B3, C3, B4, then the "Cart is empty" comments
# (H2..H6), then C4.

# --- Row 3 (Test Description / Expected Result split into 2 steps) ---
$ws.Range("B3").Value = "Add a product from the product page "
$ws.Range("C3").Value = "The selected product is added to cart successfully"

# --- Row 4 ---
$ws.Range("B4").Value = "Go to the product detail page to add product to cart of the same color and size that was previously added and verify added product information in cart"

# --- Comment column updated to "Cart is empty" for all data rows ---
$ws.Range("H2").Value = "Cart is empty"
$ws.Range("H3").Value = "Cart is empty"
$ws.Range("H4").Value = "Cart is empty"
$ws.Range("H5").Value = "Cart is empty"
$ws.Range("H6").Value = "Cart is empty"

# --- Row 4 Expected Result ---
$ws.Range("C4").Value = "The information of the added products matches the one in the cart"

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 73.5
$ws.Rows.Item(4).RowHeight = 86.25

# --- Selection ---
$ws.Range("D4").Select()
